$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.872.56'
$ws.Range('E2').Value = '  -1.20%  '
$ws.Range('D3').Value = '1.893.55'
$ws.Range('E3').Value = '  -1.21%  '
$ws.Range('D4').Value = '''1.000'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '''0.7749'
$ws.Range('E5').Value = '  -4.45%  '
$ws.Range('D6').Value = '''244.59'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').Value = '''1.000'
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').Value = '''0.3147'
$ws.Range('E8').Value = '  -3.00%  '
$ws.Range('D9').Value = '''0.07563'
$ws.Range('E9').Value = '  +4.64%  '
$ws.Range('D10').Value = '''25.50'
$ws.Range('D11').Value = '''0.08098'
$ws.Range('E11').Value = '  +0.03%  '
$ws.Range('D12').Value = '''0.7724'
$ws.Range('E12').Value = '  -2.04%  '
$ws.Range('D13').Value = '''5.505'
$ws.Range('E13').Value = '  +1.88%  '
$ws.Range('D14').Value = '1.894.69'
$ws.Range('E14').Value = '  -1.10%  '
$ws.Range('D15').Value = '''92.42'
$ws.Range('E15').Value = '  -1.42%  '
$ws.Range('D16').Value = '''6.264'
$ws.Range('E16').Value = '  +3.55%  '
$ws.Range('D17').Value = '29.966.24'
$ws.Range('E17').Value = '  -0.98%  '
$ws.Range('D18').Value = '''14.02'
$ws.Range('E18').Value = '  -1.20%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '''0.000007949'
$ws.Range('E19').Value = '  +1.36%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '''244.33'
$ws.Range('E20').Value = '  -1.96%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').Value = '''8.175'
$ws.Range('E21').Value = '  -0.72%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '''1.000'
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').Value = '2.130.24'
$ws.Range('E23').Value = '  -1.92%  '
$ws.Range('D24').Value = '''0.9993'
$ws.Range('E24').Value = '  -0.31%  '
$ws.Range('D25').Value = '''0.1566'
$ws.Range('E25').Value = '  -5.88%  '
$ws.Range('D26').Value = '''9.466'
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').Value = '''162.78'
$ws.Range('E27').Value = '  -3.18%  '
$ws.Range('D28').Value = '''18.78'
$ws.Range('E28').Value = '  -1.10%  '
$ws.Range('D29').Value = '''2.045'
$ws.Range('E29').Value = '  -5.51%  '
$ws.Range('D30').Value = '''1.441'
$ws.Range('E30').Value = '  +4.04%  '
$ws.Range('D31').Value = '''1.551'
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('D32').Value = '''4.490'
$ws.Range('E32').Value = '  +3.85%  '
$ws.Range('D33').Value = '''4.102'
$ws.Range('E33').Value = '  -0.92%  '
$ws.Range('D34').Value = '''0.05528'
$ws.Range('E34').Value = '  -5.03%  '
$ws.Range('D35').Value = '''1.259'
$ws.Range('E35').Value = '  -2.57%  '
$ws.Range('D36').Value = '''0.7592'
$ws.Range('E36').Value = '  +1.53%  '
$ws.Range('D37').Value = '''1.000'
$ws.Range('E37').Value = '  +0.44%  '
$ws.Range('D38').Value = '''2.643'
$ws.Range('E38').Value = '  -3.29%  '
$ws.Range('D39').Value = '''0.01930'
$ws.Range('E39').Value = '  -1.65%  '
$ws.Range('D40').Value = '''2.790'
$ws.Range('E40').Value = '  -1.04%  '
$ws.Range('D41').Value = '1.157.01'
$ws.Range('E41').Value = '  +13.50%  '
$ws.Range('E42').Value = '  -0.64%  '
$ws.Range('D43').Value = '''0.4437'
$ws.Range('E43').Value = '  -2.16%  '
$ws.Range('D44').Value = '''5.940'
$ws.Range('E44').Value = '  -0.53%  '
$ws.Range('D45').Value = '''0.8478'
$ws.Range('E45').Value = '  -0.55%  '
$ws.Range('D46').Value = '''0.9996'
$ws.Range('E46').Value = '  -0.15%  '
$ws.Range('E47').Value = '  -1.41%  '
$ws.Range('B48').Value = 'SynthetixNetwork'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D48').Value = '''3.146'
$ws.Range('E48').Value = '  +1.31%  '
$ws.Range('D49').Value = '''10.05'
$ws.Range('E49').Value = '  +0.84%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').Value = '''102.66'
$ws.Range('E50').Value = '  -0.65%  '
$ws.Range('D51').Value = '''7.529'
$ws.Range('E51').Value = '  -1.06%  '

# Quote-prefixing a numeric-looking string applies an implicit
# "quote prefix" cell style in Excel; reset style back to Normal so
# we only change the value (matches original unstyled cells).
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
